# Add new dimensions for base station
#
# RV1 ("3309W-1-103") is replaced with a new trim-pot part "3319W-1-103"
# that has updated cost figures. Updating the shared-string text and the
# per-unit cost cells causes all of the dependent SUM()/ratio formulas
# (column I/K/M on row 50 as well as the totals on rows 58-59) to
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 50 corresponds to component RV1 in column A.
$ws.Range("B50").Value = "3319W-1-103"
$ws.Range("H50").Value = 0.4
$ws.Range("J50").Value = 0.319
$ws.Range("L50").Value = 0.286

# Force a full recalculation so the dependent formula cells
# (I50, K50, M50, I58, K58, M58, I59, K59, M59) pick up the new values.
$excel.CalculateFull()
